$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matches original inlineStr type),
# since plain numeric-looking strings like "246.16" would otherwise be
# auto-converted to numbers by Excel on assignment.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "30.916.84"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "1.911.19"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "246.16"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D7").Value = "0.5003"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").Value = "0.2987"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.06857"
$ws.Range("E9").Value = "  +3.67%  "
$ws.Range("D10").Value = "1.907.92"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").Value = "17.02"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").Value = "0.07334"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").Value = "91.39"
$ws.Range("E13").Value = "  +6.65%  "
$ws.Range("D14").Value = "5.111"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "0.6823"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "30.892.05"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "0.000008050"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "13.31"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "2.154.14"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("D23").Value = "182.98"
$ws.Range("E23").Value = "  +34.49%  "
$ws.Range("D24").Value = "6.131"
$ws.Range("E24").Value = "  +9.51%  "
$ws.Range("D25").Value = "9.374"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Value = "154.66"
$ws.Range("E26").Value = "  +2.78%  "
$ws.Range("D27").Value = "18.78"
$ws.Range("E27").Value = "  +11.98%  "
$ws.Range("D28").Value = "1.949"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "1.394"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "4.356"
$ws.Range("E30").Value = "  +4.57%  "
$ws.Range("D31").Value = "0.09004"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  +3.20%  "
$ws.Range("D33").Value = "0.05272"
$ws.Range("E33").Value = "  +5.66%  "
$ws.Range("D34").Value = "0.7496"
$ws.Range("E34").Value = "  +6.75%  "
$ws.Range("D35").Value = "1.139"
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("D36").Value = "2.670"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "0.01943"
$ws.Range("E37").Value = "  +18.27%  "
$ws.Range("D38").Value = "2.744"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").Value = "0.9374"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "0.4419"
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("D42").Value = "106.58"
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("D43").Value = "5.862"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "7.789"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("E46").Value = "  +6.88%  "
$ws.Range("D47").Value = "0.05852"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "0.3934"
$ws.Range("E48").Value = "  +6.21%  "
$ws.Range("D49").Value = "8.625"
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("D50").Value = "33.40"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("D51").Value = "1.394"
$ws.Range("E51").Value = "  +3.90%  "

# Restore original (default) number format/style so only the cell values differ.
$priceCol.NumberFormat = "General"
$priceCol.Style = "Normal"

